# Swap the Fecha/Precio mínimo/Precio máximo/Precio promedio ponderado/Origen/
# Precio $-Kg values between row 2 and row 3 (both rows describe the same
# market+product; the weekly refresh re-paired which date goes with which
# price/origin figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $cellA = $ws.Range($col + "2")
    $cellB = $ws.Range($col + "3")

    $valA = $cellA.Value()
    $valB = $cellB.Value()

    $cellA.Value = $valB
    $cellB.Value = $valA
}
